$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill column A first (IDs), then column B (acciones), then column C
# (resultados esperados) -- mirrors the order the strings were
# originally authored in, keeping the shared-string table in sync.
$ws.Range("A23").Value = "Caso #17"
$ws.Range("A24").Value = "Caso #18"
$ws.Range("A25").Value = "Caso #19"

$ws.Range("B23").Value = "Ver reseña (admin)"
$ws.Range("B24").Value = "Editar reseña"
$ws.Range("B25").Value = "Eliminar reseña"

$ws.Range("C23").Value = "Ver reseña desde el menu de admin"
$ws.Range("C24").Value = "Editar reseña"
$ws.Range("C25").Value = "Eliminar reseña"

$ws.Range("D23").Value = "12/14/2023"
$ws.Range("D24").Value = "12/14/2023"
$ws.Range("D25").Value = "12/14/2023"

$ws.Range("E23").Value = "SI"
$ws.Range("E24").Value = "SI"
$ws.Range("E25").Value = "SI"

$ws.Range("F23").Value = "-"
$ws.Range("F24").Value = "-"
$ws.Range("F25").Value = "-"

$ws.Range("G23").Value = "OK"
$ws.Range("G24").Value = "OK"
$ws.Range("G25").Value = "OK"

# Update the sheet view: remove the frozen top-left scroll position and
# change the active selection to E23:G25
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E23:G25").Select()
